# The Power column (F) had been stored 1000x too small (e.g. Watts
# instead of milliwatts / a stray unit-conversion). Multiply every
# data row's Power value by 1000 to correct it.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 6)
    $current = $cell.Value()
    if ($current -ne $null) {
        $cell.Value = $current * 1000
    }
}

Write-Host "Rescaled Power column (F2:F$lastRow) by x1000"
